$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExecutionPlan")

# Update existing rows: mark Execute = "Yes" for rows 4 and 11
$ws.Range("D4").Value = "Yes"
$ws.Range("D11").Value = "Yes"

# Add new row 17: Home Page validation UI
$ws.Range("A17").Value = "Home Page validation UI"
$ws.Range("B17").Value = "Home Page validation UI"
$ws.Range("D17").Value = "Yes"

# Add new row 18: test script for add new person
$ws.Range("A18").Value = "test script for add new person"
$ws.Range("B18").Value = "test script for add new person"
$ws.Range("D18").Value = "Yes"
